$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the header text in B1 from "Invoice" to "Lease Agreement"
$ws.Range("B1").Value = "Lease Agreement"

# Update column widths: column B becomes wider (~18.71 chars) while others stay the same (14.7109375)
$ws.Columns.Item(2).ColumnWidth = 17.8333333333

# Give row 1 an explicit custom height
$ws.Rows.Item(1).RowHeight = 18.75

# Move the active selection from A2 to B2
$ws.Range("B2").Select()
